$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrangeHRM")

# ---- Row 1 (headers) ----
$ws.Range("F1").Value = "Message"
$ws.Range("G1").Value = "FirstName"
$ws.Range("H1").Value = "Middle Name"
$ws.Range("I1").Value = "LastName"
$ws.Range("J1").Value = "Location"
$ws.Range("K1").Value = "Nationality"
$ws.Range("L1").Value = "Marital Status"
$ws.Range("M1").Value = "Gender"
$ws.Range("N1").Value = "Region"
$ws.Range("O1").Value = "FTE"
$ws.Range("P1").Value = "Temporary Department"

# ---- Row 2 (data) ----
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "Jinny"
$ws.Range("H2").Value = "and"
$ws.Range("I2").Value = "Georgia"
$ws.Range("J2").Value = "Canadian Regional HQ"
$ws.Range("K2").Value = "Australian"
$ws.Range("L2").Value = "single"
$ws.Range("M2").Value = "Male"
$ws.Range("N2").Value = "Region-2"
$ws.Range("O2").Value = 0.5
$ws.Range("P2").Value = "Sub unit-3"

# ---- Apply center formatting to the full used range (new + existing cells) ----
$ws.Range("A1:P7").HorizontalAlignment = -4108

# ---- View / selection ----
$ws.Range("P1").Select()

Write-Output "done"
